$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 16:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 767189
$ws.Range("C4").Value = 2553
$ws.Range("D4").Value = 71281
$ws.Range("E4").Value = 655165
$ws.Range("G4").Value = 168
$ws.Range("H4").Value = 40743

# Row 28 - Arabia Saudita
$ws.Range("F28").Value = 88

# Row 38 - Noruega
$ws.Range("B38").Value = 7122
$ws.Range("C38").Value = 44
$ws.Range("E38").Value = 6919
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = 171

# Row 48 - Republica Dominicana
$ws.Range("B48").Value = 4964
$ws.Range("C48").Value = 284
$ws.Range("D48").Value = 416
$ws.Range("E48").Value = 4313
$ws.Range("F48").Value = 128
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 235

# Row 60 - Moldavia
$ws.Range("B60").Value = 2548
$ws.Range("C60").Value = 76
$ws.Range("E60").Value = 2023

# Row 68 - Uzbekistan
$ws.Range("B68").Value = 1604
$ws.Range("C68").Value = 39
$ws.Range("E68").Value = 1338
